# Update column F (dSF) values for specific rows to reflect the
# repulled/recalculated data, per the commit:
# "repull data, push all data, mean calculation"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    7  = -6
    10 = -3
    23 = 0
    26 = -2
    34 = -1
    35 = 0
    47 = -1
    48 = 2
    52 = 0
    56 = 0
    57 = 2
    64 = 5
    66 = -2
    68 = -1
    70 = -3
    71 = -5
    72 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
